# Update automàtic: dades i banners [2026-02-23 20:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Helper: assign a literal text value to a cell without Excel re-interpreting
# percentage-looking strings (e.g. "68%") as a numeric percent value.
# Strategy: temporarily force Text number format while assigning the value
# (so the "NN%" string is stored verbatim), then paste-special just the
# number format from a known General-formatted cell back onto the target —
# this resets NumberFormat to General without touching the already-stored
# text value, so the cell keeps its original style index (General/"s=3").
$fmtRef = $ws.Range("I2")
$fmtRef.Copy()
function Set-LiteralText($addr, $text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.PasteSpecial(-4122)
}

$ws.Range("E2").Value = "2026-02-23 20:18:25"
$ws.Range("E3").Value = "2026-02-23 20:18:28"
$ws.Range("E4").Value = "2026-02-23 20:18:30"
Set-LiteralText "H4" "68%"
$ws.Range("E5").Value = "2026-02-23 20:18:33"
$ws.Range("O5").Value = "4.7 °C"
$ws.Range("E6").Value = "2026-02-23 20:18:35"
$ws.Range("E7").Value = "2026-02-23 20:18:38"
$ws.Range("J7").Value = "1024.7 hPa"
$ws.Range("E8").Value = "2026-02-23 20:18:40"
$ws.Range("J8").Value = "1024.3 hPa"
$ws.Range("E9").Value = "2026-02-23 20:18:43"
Set-LiteralText "H9" "72%"
$ws.Range("O9").Value = "12.6 °C"
$ws.Range("E10").Value = "2026-02-23 20:18:45"
$ws.Range("O10").Value = "11.0 °C"
$ws.Range("E11").Value = "2026-02-23 20:18:48"
$ws.Range("E12").Value = "2026-02-23 20:18:50"
$ws.Range("E13").Value = "2026-02-23 20:18:52"
$ws.Range("E14").Value = "2026-02-23 20:18:55"
Set-LiteralText "H14" "74%"
$ws.Range("O14").Value = "12.7 °C"
$ws.Range("E15").Value = "2026-02-23 20:18:57"
$ws.Range("O15").Value = "12.7 °C"
$ws.Range("E16").Value = "2026-02-23 20:19:00"
Set-LiteralText "H16" "21%"
$ws.Range("E17").Value = "2026-02-23 20:19:02"
$ws.Range("K17").Value = "17.2 MJ/m2"
$ws.Range("E18").Value = "2026-02-23 20:19:05"
$ws.Range("J18").Value = "1024.9 hPa"
$ws.Range("O18").Value = "11.1 °C"
$ws.Range("E19").Value = "2026-02-23 20:19:08"
$ws.Range("E20").Value = "2026-02-23 20:19:10"
Set-LiteralText "H20" "38%"
$ws.Range("E21").Value = "2026-02-23 20:19:12"
$ws.Range("E22").Value = "2026-02-23 20:19:15"
$ws.Range("E23").Value = "2026-02-23 20:19:17"
Set-LiteralText "H23" "22%"
$ws.Range("E24").Value = "2026-02-23 20:19:20"
Set-LiteralText "H24" "82%"
$ws.Range("E25").Value = "2026-02-23 20:19:22"
Set-LiteralText "H25" "27%"
$ws.Range("O25").Value = "5.9 °C"
$ws.Range("E26").Value = "2026-02-23 20:19:24"
$ws.Range("J26").Value = "1023.9 hPa"
$ws.Range("O26").Value = "10.1 °C"
$ws.Range("E27").Value = "2026-02-23 20:19:27"
$ws.Range("E28").Value = "2026-02-23 20:19:29"
$ws.Range("E29").Value = "2026-02-23 20:19:32"
Set-LiteralText "H29" "82%"
$ws.Range("O29").Value = "10.9 °C"
$ws.Range("E30").Value = "2026-02-23 20:19:34"
$ws.Range("J30").Value = "1024.5 hPa"
$ws.Range("E31").Value = "2026-02-23 20:19:37"
$ws.Range("E32").Value = "2026-02-23 20:19:39"
$ws.Range("O32").Value = "7.9 °C"
$ws.Range("E33").Value = "2026-02-23 20:19:41"
$ws.Range("E34").Value = "2026-02-23 20:19:44"
Set-LiteralText "H34" "42%"
$ws.Range("E35").Value = "2026-02-23 20:19:46"
$ws.Range("N35").Value = "6.8 °C 19:54 TU"
$ws.Range("O35").Value = "12.5 °C"
$ws.Range("E36").Value = "2026-02-23 20:19:48"
Set-LiteralText "H36" "73%"
$ws.Range("J36").Value = "1024.6 hPa"
$ws.Range("E37").Value = "2026-02-23 20:19:51"
$ws.Range("O37").Value = "9.3 °C"
$ws.Range("E38").Value = "2026-02-23 20:19:53"
Set-LiteralText "H38" "64%"
$ws.Range("E39").Value = "2026-02-23 20:19:55"
$ws.Range("E40").Value = "2026-02-23 20:19:58"
Set-LiteralText "H40" "61%"
$ws.Range("E41").Value = "2026-02-23 20:20:00"
Set-LiteralText "H41" "72%"
$ws.Range("O41").Value = "12.2 °C"
$ws.Range("E42").Value = "2026-02-23 20:20:02"
$ws.Range("E43").Value = "2026-02-23 20:20:05"
$ws.Range("E44").Value = "2026-02-23 20:20:07"
$ws.Range("N44").Value = "0.2 °C 19:45 TU"
$ws.Range("O44").Value = "3.2 °C"
$ws.Range("E45").Value = "2026-02-23 20:20:10"
$ws.Range("O45").Value = "8.5 °C"
$ws.Range("E46").Value = "2026-02-23 20:20:12"

$excel.CutCopyMode = 0

